$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update grade data (Tema_4 = column F, Tema_5 = column G; a couple of
#     Tema_3 / column E corrections) for newly graded submissions ---
$ws.Range("F2").Value = 60
$ws.Range("G2").Value = 54
$ws.Range("F3").Value = 59
$ws.Range("G3").Value = 47
$ws.Range("F4").Value = 46
$ws.Range("G4").Value = 55
$ws.Range("E6").Value = 46
$ws.Range("G6").Value = 51
$ws.Range("F7").Value = 40
$ws.Range("G7").Value = 46
$ws.Range("G8").Value = 53
$ws.Range("G10").Value = 56
$ws.Range("F12").Value = 32
$ws.Range("G12").Value = 40
$ws.Range("F13").Value = 54
$ws.Range("G13").Value = 55
$ws.Range("F14").Value = 57
$ws.Range("G14").Value = 55
$ws.Range("F15").Value = 25
$ws.Range("G15").Value = 48
$ws.Range("G17").Value = 41
$ws.Range("F18").Value = 58
$ws.Range("G18").Value = 53
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 47
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 54
$ws.Range("F23").Value = 55
$ws.Range("G23").Value = 48
$ws.Range("G26").Value = 44
$ws.Range("F27").Value = 52
$ws.Range("G27").Value = 52
$ws.Range("F28").Value = 52
$ws.Range("G28").Value = 45
$ws.Range("F29").Value = 58
$ws.Range("G31").Value = 52
$ws.Range("G33").Value = 52
$ws.Range("F34").Value = 50
$ws.Range("G34").Value = 51
$ws.Range("F35").Value = 52
$ws.Range("G35").Value = 52
$ws.Range("F38").Value = 43
$ws.Range("G38").Value = 52
$ws.Range("G39").Value = 39
$ws.Range("F40").Value = 64
$ws.Range("G40").Value = 53
$ws.Range("G41").Value = 53
$ws.Range("F42").Value = 46
$ws.Range("G42").Value = 27
$ws.Range("G44").Value = 36
$ws.Range("F45").Value = 51
$ws.Range("G46").Value = 53
$ws.Range("F47").Value = 60
$ws.Range("G47").Value = 53
$ws.Range("F48").Value = 28
$ws.Range("F50").Value = 67
$ws.Range("G50").Value = 51
$ws.Range("F52").Value = 65
$ws.Range("G52").Value = 52
$ws.Range("F53").Value = 67
$ws.Range("G53").Value = 52
$ws.Range("F55").Value = 45
$ws.Range("F56").Value = 48
$ws.Range("G56").Value = 52
$ws.Range("F60").Value = 63
$ws.Range("G60").Value = 52
$ws.Range("F61").Value = 58
$ws.Range("G61").Value = 52
$ws.Range("F63").Value = 66
$ws.Range("G63").Value = 52
$ws.Range("F64").Value = 65
$ws.Range("G64").Value = 52
$ws.Range("F65").Value = 65
$ws.Range("G65").Value = 52
$ws.Range("G66").Value = 52
$ws.Range("F67").Value = 64
$ws.Range("G67").Value = 52
$ws.Range("F68").Value = 64
$ws.Range("G68").Value = 52
$ws.Range("F71").Value = 66
$ws.Range("G71").Value = 52
$ws.Range("F72").Value = 58
$ws.Range("F73").Value = 64
$ws.Range("G73").Value = 52
$ws.Range("G75").Value = 53
$ws.Range("F76").Value = 68
$ws.Range("G76").Value = 57
$ws.Range("F77").Value = 43
$ws.Range("G77").Value = 47
$ws.Range("F78").Value = 52
$ws.Range("G78").Value = 53
$ws.Range("F79").Value = 62
$ws.Range("G79").Value = 53
$ws.Range("F80").Value = 28
$ws.Range("G80").Value = 52
$ws.Range("F81").Value = 57
$ws.Range("G81").Value = 49
$ws.Range("F82").Value = 49
$ws.Range("G82").Value = 53
$ws.Range("F84").Value = 59
$ws.Range("G84").Value = 53
$ws.Range("F85").Value = 54
$ws.Range("G85").Value = 53
$ws.Range("G86").Value = 53
$ws.Range("F88").Value = 23
$ws.Range("F89").Value = 62
$ws.Range("G89").Value = 54
$ws.Range("F90").Value = 54
$ws.Range("G90").Value = 53
$ws.Range("F91").Value = 49
$ws.Range("F92").Value = 65
$ws.Range("G92").Value = 52
$ws.Range("F93").Value = 64
$ws.Range("G93").Value = 54
$ws.Range("F94").Value = 34
$ws.Range("G94").Value = 53
$ws.Range("F95").Value = 52
$ws.Range("G95").Value = 53
$ws.Range("F96").Value = 62
$ws.Range("G96").Value = 53
$ws.Range("F97").Value = 54
$ws.Range("G97").Value = 53
$ws.Range("F98").Value = 63
$ws.Range("G98").Value = 53
$ws.Range("F99").Value = 18
$ws.Range("G99").Value = 52
$ws.Range("F101").Value = 65
$ws.Range("G101").Value = 53
$ws.Range("F102").Value = 60
$ws.Range("G102").Value = 53
$ws.Range("G103").Value = 53
$ws.Range("F104").Value = 28
$ws.Range("G104").Value = 51
$ws.Range("F105").Value = 26
$ws.Range("G105").Value = 37
$ws.Range("F106").Value = 28
$ws.Range("G106").Value = 48

# --- Unhide all rows and clear the autofilter criteria so every student row
#     (previously filtered to show only group 244) is visible again ---
$ws.ShowAllData()

# --- Restore the active selection recorded in the saved workbook ---
$ws.Range("C24").Select()
